$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Range("H18").Value = 502.73334
$ws.Range("I18").Value = 474.1
$ws.Range("J18").Value = 560
$ws.Range("K18").Value = 474.1
$ws.Range("L18").Value = 560
$ws.Range("M18").Value = -190.1
$ws.Range("N18").Value = -1128

# Row 28
$ws.Range("H28").Value = 2623.75
$ws.Range("I28").Value = 2498.3333
$ws.Range("J28").Value = 3000
$ws.Range("K28").Value = 2498.3333
$ws.Range("L28").Value = 3000
$ws.Range("M28").Value = -2013.3333
$ws.Range("N28").Value = -3970

# Row 86
$ws.Range("H86").Value = 1208.1666
$ws.Range("I86").Value = 999
$ws.Range("J86").Value = 1250
$ws.Range("K86").Value = 999
$ws.Range("L86").Value = 1250
$ws.Range("M86").Value = 124
$ws.Range("N86").Value = -3496

# Row 89
$ws.Range("H89").Value = 1208.1666
$ws.Range("I89").Value = 999
$ws.Range("J89").Value = 1250
$ws.Range("K89").Value = 4995
$ws.Range("L89").Value = 6250
$ws.Range("M89").Value = 621
$ws.Range("N89").Value = -17482

# Row 113
$ws.Range("H113").Value = 8144.5
$ws.Range("I113").Value = 6081.273
$ws.Range("K113").Value = 6081.273
$ws.Range("M113").Value = -2827.273

# Row 137
$ws.Range("H137").Value = 1375
$ws.Range("I137").Value = 1250
$ws.Range("J137").Value = 1500
$ws.Range("K137").Value = 3750
$ws.Range("L137").Value = 4500
$ws.Range("M137").Value = -1200
$ws.Range("N137").Value = -9600

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1057.7142
$ws.Range("I2").Value = 1057.7142
$ws.Range("K2").Value = 1057.7142
$ws.Range("M2").Value = -944.7141999999999

# Row 3
$ws.Range("H3").Value = 806.6
$ws.Range("I3").Value = 496.33334
$ws.Range("J3").Value = 1272
$ws.Range("K3").Value = 496.33334
$ws.Range("L3").Value = 1272
$ws.Range("M3").Value = -381.33334
$ws.Range("N3").Value = -1502

# Row 35
$ws.Range("H35").Value = 2000
$ws.Range("I35").Value = 2000
$ws.Range("K35").Value = 2000
$ws.Range("M35").Value = -1594

# Row 41
$ws.Range("H41").Value = 7241
$ws.Range("J41").Value = 34000
$ws.Range("L41").Value = 34000
$ws.Range("N41").Value = -34828

# Row 61
$ws.Range("H61").Value = 1925
$ws.Range("I61").Value = 1925
$ws.Range("K61").Value = 1925
$ws.Range("M61").Value = -1713

# Row 97
$ws.Range("H97").Value = 1766.4348
$ws.Range("J97").Value = 2042.7142
$ws.Range("L97").Value = 2042.7142
$ws.Range("N97").Value = -3034.7142

# Row 116
$ws.Range("H116").Value = 1057.7142
$ws.Range("I116").Value = 1057.7142
$ws.Range("K116").Value = 1057.7142
$ws.Range("M116").Value = 1236.2858

# Row 136
$ws.Range("H136").Value = 1925
$ws.Range("I136").Value = 1925
$ws.Range("K136").Value = 5775
$ws.Range("M136").Value = -3225

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1057.7142
$ws.Range("I3").Value = 1057.7142
$ws.Range("K3").Value = 1057.7142
$ws.Range("M3").Value = -943.7141999999999

# Row 36
$ws.Range("H36").Value = 3999
$ws.Range("I36").Value = 3999
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 3999
$ws.Range("L36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -3465

# Row 54
$ws.Range("H54").Value = 3915
$ws.Range("I54").Value = 3915
$ws.Range("K54").Value = 3915
$ws.Range("M54").Value = -3431

$ws = $wb.Worksheets.Item("CRP")
# Row 107
$ws.Range("H107").Value = 582.9231
$ws.Range("I107").Value = 397.9
$ws.Range("J107").Value = 1199.6666
$ws.Range("K107").Value = 397.9
$ws.Range("L107").Value = 1199.6666
$ws.Range("M107").Value = 1522.1
$ws.Range("N107").Value = -5039.6666

# Row 122
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").ClearContents()
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = 0

$ws = $wb.Worksheets.Item("CUL")
# Row 117
$ws.Range("H117").Value = 1767.875
$ws.Range("J117").Value = 3198.75
$ws.Range("L117").Value = 9596.25
$ws.Range("N117").Value = -16480.25

$ws = $wb.Worksheets.Item("GSM")
# Row 41
$ws.Range("H41").Value = 12500
$ws.Range("I41").Value = 12500
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 12500
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -12145

# Row 80
$ws.Range("H80").Value = 2751.6667
$ws.Range("I80").Value = 2315
$ws.Range("J80").Value = 3625
$ws.Range("K80").Value = 2315
$ws.Range("L80").Value = 3625
$ws.Range("M80").Value = -1317
$ws.Range("N80").Value = -5621

# Row 83
$ws.Range("H83").Value = 2751.6667
$ws.Range("I83").Value = 2315
$ws.Range("J83").Value = 3625
$ws.Range("K83").Value = 11575
$ws.Range("L83").Value = 18125
$ws.Range("M83").Value = -6583
$ws.Range("N83").Value = -28109

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 3427.1428
$ws.Range("I46").Value = 997.5
$ws.Range("K46").Value = 997.5
$ws.Range("M46").Value = -809.5

# Row 82
$ws.Range("H82").Value = 1582
$ws.Range("I82").Value = 1462.375
$ws.Range("J82").Value = 1821.25
$ws.Range("K82").Value = 1462.375
$ws.Range("L82").Value = 1821.25
$ws.Range("M82").Value = -1101.375
$ws.Range("N82").Value = -2543.25

# Row 85
$ws.Range("H85").Value = 1582
$ws.Range("I85").Value = 1462.375
$ws.Range("J85").Value = 1821.25
$ws.Range("K85").Value = 1462.375
$ws.Range("L85").Value = 1821.25
$ws.Range("M85").Value = -214.375
$ws.Range("N85").Value = -4317.25

# Row 100
$ws.Range("H100").Value = 2375.75
$ws.Range("I100").Value = 2375.75
$ws.Range("K100").Value = 2375.75
$ws.Range("M100").Value = -1834.75

$ws = $wb.Worksheets.Item("WVR")
# Row 23
$ws.Range("H23").Value = 2236.6667
$ws.Range("J23").Value = 5950
$ws.Range("L23").Value = 5950
$ws.Range("N23").Value = -6408

# Row 51
$ws.Range("H51").Value = 17500
$ws.Range("I51").Value = 5000
$ws.Range("K51").Value = 5000
$ws.Range("M51").Value = -4490

# Row 81
$ws.Range("H81").Value = 4099.2856
$ws.Range("I81").Value = 4115.5
$ws.Range("K81").Value = 8231
$ws.Range("M81").Value = -7170

# Row 84
$ws.Range("H84").Value = 4099.2856
$ws.Range("I84").Value = 4115.5
$ws.Range("K84").Value = 41155
$ws.Range("M84").Value = -35851

# Row 100
$ws.Range("H100").Value = 731.7778
$ws.Range("I100").Value = 698.375
$ws.Range("K100").Value = 1396.75
$ws.Range("M100").Value = -855.75

# Row 107
$ws.Range("H107").Value = 1982
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()
